# COREESG_holdings.xlsx update
# - bump the "as of" date in the confidential disclaimer (A10) from 2021-05-21 to 2021-05-24
# - refresh the Weight / Percent Change figures in D2:E7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected (legacy password "D382"); unlock it so the cells can be written.
$ws.Unprotect("D382")

# Updated disclaimer text (only the date changed).
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-24 for illustrative purposes only and are subject to change."
# Undo any row auto-fit side effect from writing the multi-line string so row 10 keeps its default height.
$ws.Rows(10).AutoFit()

# Updated Weight (D) / Percent Change (E) values.
$ws.Range("D2").Value = 0.2451526361615041
$ws.Range("E2").Value = 0.01662747732616743

$ws.Range("D3").Value = 0.5007422904811412
$ws.Range("E3").Value = 0.004207204838285294

$ws.Range("D4").Value = 0.09624569995190031
$ws.Range("E4").Value = 0.009699129057798794

$ws.Range("D5").Value = 0.1015610982267165
$ws.Range("E5").Value = 0.005737704918032938

$ws.Range("D6").Value = 0.05629827517873801
$ws.Range("E6").Value = 0.005246350364963348

$ws.Range("E7").Value = 0.007994582840870512

# Restore sheet protection with the original password.
$ws.Protect("D382")
